# Updates odds values on Sheet1 to reflect the latest FlashScore data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of cell address -> new value
$updates = @{
    # Row 3
    "G3"  = 2
    "H3"  = 3
    "I3"  = 3.95
    "J3"  = 2.55
    "L3"  = 4.3
    "W3"  = 6.5
    "X3"  = 9.25
    "Y3"  = 8.25
    "Z3"  = 18.5
    "AA3" = 17
    "AB3" = 29
    "AH3" = 10.25
    "AI3" = 22
    "AJ3" = 13
    "AK3" = 65
    "AL3" = 40
    "AN3" = 3.85
    "AO3" = 10
    "AP3" = 17.5
    "AQ3" = 37
    "AR3" = 65
    "AW3" = 5.8
    "AX3" = 22
    "AY3" = 26
    "AZ3" = 120
    "BA3" = 150

    # Row 5
    "Q5" = 2.03
    "R5" = 1.83

    # Row 7
    "G7"  = 3.3
    "I7"  = 2.15
    "J7"  = 4
    "L7"  = 3
    "N7"  = 7.5
    "Q7"  = 2.3
    "R7"  = 1.6
    "W7"  = 8.5
    "X7"  = 15
    "Y7"  = 12
    "AA7" = 29
    "AF7" = 51
    "AI7" = 9.5
    "AK7" = 21
    "AL7" = 21
    "AN7" = 5
    "AQ7" = 67
    "AS7" = 251
    "AX7" = 13

    # Row 9
    "O9" = 1.3
    "P9" = 3.4
    "Q9" = 2
    "R9" = 1.85
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
